$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 44372
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 50
$ws.Range("K3").Value = 6000
$ws.Range("L3").Value = 7000
$ws.Range("M3").Value = 6400
$ws.Range("O3").Value = "Región del Maule"
$ws.Range("P3").Value = 400

# Row 4
$ws.Range("D4").Value = 44371
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 7000
$ws.Range("L4").Value = 7000
$ws.Range("M4").Value = 7000
$ws.Range("O4").Value = "Región Metropolitana"
$ws.Range("P4").Value = 438

# Row 5
$ws.Range("D5").Value = 44308
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 75
$ws.Range("K5").Value = 5000
$ws.Range("L5").Value = 5000
$ws.Range("M5").Value = 5000
$ws.Range("O5").Value = "Región del Maule"
$ws.Range("P5").Value = 312

# Row 6
$ws.Range("D6").Value = 44369
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 60
$ws.Range("K6").Value = 7000
$ws.Range("L6").Value = 7000
$ws.Range("M6").Value = 7000
$ws.Range("O6").Value = "Región Metropolitana"
$ws.Range("P6").Value = 438

# Row 7
$ws.Range("D7").Value = 44314
$ws.Range("I7").Value = "Segunda"
$ws.Range("J7").Value = 20
$ws.Range("K7").Value = 5000
$ws.Range("L7").Value = 5000
$ws.Range("M7").Value = 5000
$ws.Range("O7").Value = "Región del Maule"
$ws.Range("P7").Value = 312

# Row 8
$ws.Range("D8").Value = 44354
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 100
$ws.Range("K8").Value = 8000
$ws.Range("L8").Value = 9000
$ws.Range("M8").Value = 8500
$ws.Range("O8").Value = "Región Metropolitana"
$ws.Range("P8").Value = 531

# Row 9
$ws.Range("D9").Value = 44354
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 80
$ws.Range("K9").Value = 9000
$ws.Range("L9").Value = 9000
$ws.Range("M9").Value = 9000
$ws.Range("O9").Value = "Región del Maule"
$ws.Range("P9").Value = 562

# Row 10
$ws.Range("D10").Value = 44467
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 40
$ws.Range("K10").Value = 7000
$ws.Range("L10").Value = 7000
$ws.Range("M10").Value = 7000
$ws.Range("O10").Value = "Región del Maule"
$ws.Range("P10").Value = 438

# Row 11
$ws.Range("D11").Value = 44305
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 35
$ws.Range("K11").Value = 7000
$ws.Range("L11").Value = 7000
$ws.Range("M11").Value = 7000
$ws.Range("O11").Value = "Región del Maule"
$ws.Range("P11").Value = 438

# Row 12
$ws.Range("D12").Value = 44398
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 80
$ws.Range("K12").Value = 7000
$ws.Range("L12").Value = 7000
$ws.Range("M12").Value = 7000
$ws.Range("O12").Value = "Región Metropolitana"
$ws.Range("P12").Value = 438

# Row 13
$ws.Range("D13").Value = 44392
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 95
$ws.Range("K13").Value = 7000
$ws.Range("L13").Value = 7000
$ws.Range("M13").Value = 7000
$ws.Range("O13").Value = "Región del Maule"
$ws.Range("P13").Value = 438

# Row 14
$ws.Range("D14").Value = 44386
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 40
$ws.Range("K14").Value = 7000
$ws.Range("L14").Value = 7000
$ws.Range("M14").Value = 7000
$ws.Range("O14").Value = "Región del Maule"
$ws.Range("P14").Value = 438

# Row 15
$ws.Range("D15").Value = 44313
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 20
$ws.Range("K15").Value = 7000
$ws.Range("L15").Value = 7000
$ws.Range("M15").Value = 7000
$ws.Range("O15").Value = "Región del Maule"
$ws.Range("P15").Value = 438

# Row 16
$ws.Range("D16").Value = 44348
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 35
$ws.Range("K16").Value = 7000
$ws.Range("L16").Value = 7000
$ws.Range("M16").Value = 7000
$ws.Range("O16").Value = "Región del Maule"
$ws.Range("P16").Value = 438

# Row 17
$ws.Range("D17").Value = 44420
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 45
$ws.Range("K17").Value = 8000
$ws.Range("L17").Value = 8000
$ws.Range("M17").Value = 8000
$ws.Range("O17").Value = "Región Metropolitana"
$ws.Range("P17").Value = 500

# Row 18
$ws.Range("D18").Value = 44315
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 40
$ws.Range("K18").Value = 7000
$ws.Range("L18").Value = 7000
$ws.Range("M18").Value = 7000
$ws.Range("O18").Value = "Región del Maule"
$ws.Range("P18").Value = 438

# Row 19
$ws.Range("D19").Value = 44403
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 35
$ws.Range("K19").Value = 5000
$ws.Range("L19").Value = 5000
$ws.Range("M19").Value = 5000
$ws.Range("O19").Value = "Región Metropolitana"
$ws.Range("P19").Value = 312

# Row 20
$ws.Range("D20").Value = 44362
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 25
$ws.Range("K20").Value = 8000
$ws.Range("L20").Value = 8000
$ws.Range("M20").Value = 8000
$ws.Range("O20").Value = "Región Metropolitana"
$ws.Range("P20").Value = 500

# Row 21
$ws.Range("D21").Value = 44397
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 40
$ws.Range("K21").Value = 8000
$ws.Range("L21").Value = 8000
$ws.Range("M21").Value = 8000
$ws.Range("O21").Value = "Región Metropolitana"
$ws.Range("P21").Value = 500

# Row 22
$ws.Range("D22").Value = 44355
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 30
$ws.Range("K22").Value = 8000
$ws.Range("L22").Value = 8000
$ws.Range("M22").Value = 8000
$ws.Range("O22").Value = "Región Metropolitana"
$ws.Range("P22").Value = 500

# Row 23
$ws.Range("D23").Value = 44396
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 80
$ws.Range("K23").Value = 7000
$ws.Range("L23").Value = 7000
$ws.Range("M23").Value = 7000
$ws.Range("O23").Value = "Región Metropolitana"
$ws.Range("P23").Value = 438

# Row 24
$ws.Range("D24").Value = 44385
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 100
$ws.Range("K24").Value = 7000
$ws.Range("L24").Value = 7000
$ws.Range("M24").Value = 7000
$ws.Range("O24").Value = "Región del Maule"
$ws.Range("P24").Value = 438

# Row 25
$ws.Range("D25").Value = 44312
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 40
$ws.Range("K25").Value = 7000
$ws.Range("L25").Value = 7000
$ws.Range("M25").Value = 7000
$ws.Range("O25").Value = "Región del Maule"
$ws.Range("P25").Value = 438

# Row 26
$ws.Range("D26").Value = 44399
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 80
$ws.Range("K26").Value = 7000
$ws.Range("L26").Value = 7000
$ws.Range("M26").Value = 7000
$ws.Range("O26").Value = "Región Metropolitana"
$ws.Range("P26").Value = 438
